# Update 'want to go' counts (column F) on sheet 展览 (sheet1), rows 2-25
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$sheet1F = @(1674, 9109, 113, 507, 706, 1374, 202, 58, 95, 5907, 60, 387, 103, 4436, 14, 165, 1148, 29, 338, 26, 257, 15, 2753, 127)
for ($idx = 0; $idx -lt $sheet1F.Count; $idx++) {
    $ws1.Cells.Item($idx + 2, 6).Value = $sheet1F[$idx]
}

# Update 'want to go' counts (column F) on sheet 演出 (sheet2), rows 2-3
$ws2 = $wb.Worksheets.Item("演出")
$sheet2F = @(34, 40)
for ($idx = 0; $idx -lt $sheet2F.Count; $idx++) {
    $ws2.Cells.Item($idx + 2, 6).Value = $sheet2F[$idx]
}

# Rebuild sheet 全部类型 (sheet4): merge of 展览 + 演出 rows, sorted by date
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Rows = @(
    @{B="2024-06-22"; C="合肥·Look Look动漫嘉年华"; D="新站区东方大道288号 少荃体育中心"; E="2024.06.22 10:00-06.23 17:30"; F=1674; G=58; H="https://show.bilibili.com/platform/detail.html?id=82311"; I="//i1.hdslb.com/bfs/openplatform/202406/1UVGJ3G01718620439056.jpeg"}
    @{B="2024-06-22"; C="合肥·城市动漫节"; D="包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心"; E="2024.06.22 10:00-06.23 16:30"; F=9109; G=70; H="https://show.bilibili.com/platform/detail.html?id=85000"; I="//i2.hdslb.com/bfs/openplatform/202405/BQV7zeWg1716290459878.jpeg"}
    @{B="2024-06-22"; C="合肥·城市动漫节·触手猫兽漫联合专区"; D="包河经济开发区大连路与园博大道交口骆岗中央公园园博小镇一期S6区1号楼 大机库演艺中心"; E="2024.06.22 10:00-06.23 16:30"; F=113; G=89; H="https://show.bilibili.com/platform/detail.html?id=86419"; I="//i2.hdslb.com/bfs/openplatform/202405/esdErBTC1716799359305.jpeg"}
    @{B="2024-06-29"; C="合肥·东方LiveParty×安徽THO4·隙间皖韵之梦"; D="国祯广场B-1楼 背影骑士LIVEHOUSE"; E="2024.06.29 18:00-06.29 21:00"; F=34; G=249; H="https://show.bilibili.com/platform/detail.html?id=87669"; I="//i2.hdslb.com/bfs/openplatform/202406/tzVYkFc21718180314059.jpeg"}
    @{B="2024-06-29"; C="合肥·原神X星铁Only"; D="金寨路与天堂窄路交叉口 梵木艺术中心"; E="2024.06.29 10:00-06.29 17:00"; F=507; G="不可售"; H="https://show.bilibili.com/platform/detail.html?id=86406"; I="//i1.hdslb.com/bfs/openplatform/202405/r3c5IueN1716820859877.jpeg"}
    @{B="2024-06-29"; C="合肥·穿梭次元动漫嘉年华"; D="金寨路与天堂窄路交叉口 梵木艺术中心"; E="2024.06.29 10:00-06.29 17:00"; F=706; G=60; H="https://show.bilibili.com/platform/detail.html?id=86428"; I="//i0.hdslb.com/bfs/openplatform/202405/gFcsiZHY1716820470513.jpeg"}
    @{B="2024-06-30"; C="合肥·第1.5届星芒动漫嘉年华"; D="山西路与太原路交叉口 挥动体育"; E="2024.06.30 09:30-06.30 17:30"; F=1374; G=60; H="https://show.bilibili.com/platform/detail.html?id=85213"; I="//i1.hdslb.com/bfs/openplatform/202405/v40vLtJl1715073148563.jpeg"}
    @{B="2024-06-30"; C="安徽·THO4·隙间皖韵之梦"; D="北二环与新蚌埠路交汇处 蓝金湾大酒店"; E="2024.06.30 10:00-06.30 17:00"; F=202; G=65; H="https://show.bilibili.com/platform/detail.html?id=85119"; I="//i2.hdslb.com/bfs/openplatform/202405/kuuarwvJ1714932457216.jpeg"}
    @{B="2024-07-06"; C="合肥·次元日记动漫游戏嘉年华"; D="徽州大道5558号(徽州大道与紫云路交口) 合肥方圆荟(滨湖店)"; E="2024.07.06 10:00-07.06 17:00"; F=58; G=45; H="https://show.bilibili.com/platform/detail.html?id=87201"; I="//i0.hdslb.com/bfs/openplatform/202406/BhvxoidA1717762410463.jpeg"}
    @{B="2024-07-06"; C="合肥·首届AS运动番Only"; D="昭潭路名都花园西南侧约150米 职工体育活动中心"; E="2024.07.06 09:00-07.06 17:00"; F=95; G=68; H="https://show.bilibili.com/platform/detail.html?id=87384"; I="//i1.hdslb.com/bfs/openplatform/202406/suoZa5Ha1717727447336.jpeg"}
    @{B="2024-07-13"; C="合肥·星光次元动漫文化节"; D="凤淮路与公园路交叉口南行50米路西 庐阳区全民健身中心"; E="2024.07.13 09:00-07.14 16:30"; F=5907; G=50; H="https://show.bilibili.com/platform/detail.html?id=86275"; I="//i1.hdslb.com/bfs/openplatform/202405/CnTO3a4Y1716539865739.jpeg"}
    @{B="2024-07-14"; C="合肥·第一届宅舞比赛漫展-CF01"; D="长江中路98号 合肥银泰君亭酒店"; E="2024.07.14 10:00-07.14 16:00"; F=60; G=59; H="https://show.bilibili.com/platform/detail.html?id=85544"; I="//i2.hdslb.com/bfs/openplatform/202405/sDlSCXnC1715316021248.jpeg"}
    @{B="2024-07-20"; C="合肥·W·A首届童年怀旧only"; D="铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)"; E="2024.07.20 09:30-07.20 17:00"; F=387; G=78; H="https://show.bilibili.com/platform/detail.html?id=84794"; I="//i2.hdslb.com/bfs/openplatform/202404/Ie0KTNEr1713951888990.png"}
    @{B="2024-07-20"; C="合肥·lovelive only"; D="莲花路与繁华大道交叉口258号 百乐门大剧院"; E="2024.07.20 11:00-07.20 19:30"; F=103; G=68; H="https://show.bilibili.com/platform/detail.html?id=87058"; I="//i2.hdslb.com/bfs/openplatform/202406/M6rzTe6y1717600620230.jpeg"}
    @{B="2024-07-20"; C="合肥·第十四届次元之门动漫游戏博览会"; D="南京路与庐州大道交汇处 合肥滨湖国际会展中心"; E="2024.07.20 10:00-07.21 17:00"; F=4436; G=68; H="https://show.bilibili.com/platform/detail.html?id=85336"; I="//i2.hdslb.com/bfs/openplatform/202405/Bu6iQPJ01715161445356.jpeg"}
    @{B="2024-07-20"; C="合肥·第十四届次元之门动漫游戏博览会一唐雅菁&一口井专场票"; D="南京路与庐州大道交汇处 合肥滨湖国际会展中心"; E="2024.07.20 10:00-07.20 15:00"; F=14; G=128; H="https://show.bilibili.com/platform/detail.html?id=87418"; I="//i1.hdslb.com/bfs/openplatform/202406/V5NW9yRo1718548308064.jpeg"}
    @{B="2024-07-20"; C="安徽·赛马娘Only 2.0"; D="文忠路1865号 赫拉诺言艺术中心"; E="2024.07.20 09:00-07.20 17:00"; F=165; G=78; H="https://show.bilibili.com/platform/detail.html?id=84539"; I="//i1.hdslb.com/bfs/openplatform/202405/ibcY9Edj1715235810905.jpeg"}
    @{B="2024-07-21"; C="合肥·首届Gumi同人展"; D="新站区东方大道288号 少荃体育中心"; E="2024.07.21 09:30-07.21 17:00"; F=1148; G=68; H="https://show.bilibili.com/platform/detail.html?id=86573"; I="//i0.hdslb.com/bfs/openplatform/202405/DsvnHgmP1717038341915.jpeg"}
    @{B="2024-07-27"; C="合肥·灵能百分百ONLY2.0"; D="铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)"; E="2024.07.27 10:00-07.27 17:00"; F=29; G=75; H="https://show.bilibili.com/platform/detail.html?id=87497"; I="//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg"}
    @{B="2024-07-27"; C="安徽·MAX特摄only展"; D="桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"; E="2024.07.27 09:30-07.27 18:00"; F=338; G=50; H="https://show.bilibili.com/platform/detail.html?id=83684"; I="//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg"}
    @{B="2024-07-27"; C="庐江·夏日游嘉年华"; D="白山路东150米 庐江体育馆"; E="2024.07.27 09:00-07.28 17:00"; F=26; G=60; H="https://show.bilibili.com/platform/detail.html?id=87569"; I="//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg"}
    @{B="2024-07-28"; C="合肥·咒术回战only"; D="清河路19号 依立腾工业园区"; E="2024.07.28 09:30-07.28 17:30"; F=257; G=60; H="https://show.bilibili.com/platform/detail.html?id=86520"; I="//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png"}
    @{B="2024-07-28"; C="合肥·第二届TH动漫游戏嘉年华"; D="田埠西路199号 吉祥如意宴会楼蜀山店"; E="2024.07.28 09:30-07.28 17:00"; F=15; G=55; H="https://show.bilibili.com/platform/detail.html?id=87447"; I="//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png"}
    @{B="2024-08-03"; C="合肥·第七届环形宇宙动漫游戏嘉年华"; D="南京路与庐州大道交汇处 合肥滨湖国际会展中心"; E="2024.08.03 09:30-08.04 17:00"; F=2753; G=49; H="https://show.bilibili.com/platform/detail.html?id=84767"; I="//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg"}
    @{B="2024-08-03"; C="合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会"; D="徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"; E="2024.08.03 19:30-08.03 21:00"; F=40; G=80; H="https://show.bilibili.com/platform/detail.html?id=83556"; I="//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg"}
    @{B="2024-08-17"; C="合肥·银魂主题派对only2.0"; D="长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)"; E="2024.08.17 13:00-08.17 18:00"; F=127; G=128; H="https://show.bilibili.com/platform/detail.html?id=87173"; I="//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png"}
)

$r = 2
foreach ($row in $sheet4Rows) {
    $ws4.Cells.Item($r, 1).Value = $r - 1
    $ws4.Cells.Item($r, 2).Value = "'" + $row.B
    $ws4.Cells.Item($r, 3).Value = $row.C
    $ws4.Cells.Item($r, 4).Value = $row.D
    $ws4.Cells.Item($r, 5).Value = $row.E
    $ws4.Cells.Item($r, 6).Value = $row.F
    $ws4.Cells.Item($r, 7).Value = $row.G
    $ws4.Cells.Item($r, 8).Value = $row.H
    $ws4.Cells.Item($r, 9).Value = $row.I
    $r++
}

Write-Host "done"